# ---------------------------------------------------------------------------
# Applies the "additional scraping" edit:
#   1. Adds a new "Player Info" sheet (before "ODI Batting") with player bio.
#   2. Renames "MATCH_CARD_LINK" -> "MATCH_CODE" on "ODI Batting" and
#      "ODI Bowling", replacing the full scorecard URL with just the
#      numeric match code, and removing a handful of spurious empty cells.
#   3. Adds a new "ODI Batting Extra" sheet (after "ODI Bowling") with
#      additional per-innings batting detail.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before the current first sheet.
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

$playerInfo.Cells.Item(2,1).Value = "4690"
$playerInfo.Cells.Item(2,2).Value = "Faheem Ashraf"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Medium"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet - rename MATCH_CARD_LINK -> MATCH_CODE (col D),
#    keep just the numeric match code, drop a few stray empty B cells.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Cells.Item(1,4).Value = "MATCH_CODE"

$battingCodes = @(
    "4044","4082","4087","4103","4104","4105","4110","4114","4172","4174",
    "4176","4177","4194","4197","4225","4227","4237","4238","4273","4274",
    "4287","4292","4294","4432","4433","4458","4459","4460","4472","4473",
    "4476"
)
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $odiBatting.Cells.Item($i + 2, 4).Value = $battingCodes[$i]
}

$emptyBattingRows = @(3,11,12,13,14,16,18,22,26)
foreach ($r in $emptyBattingRows) {
    $odiBatting.Cells.Item($r, 2).Value = ""
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet - rename MATCH_CARD_LINK -> MATCH_CODE (col B),
#    keep just the numeric match code.
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Cells.Item(1,2).Value = "MATCH_CODE"

$bowlingCodes = @(
    "4044","4082","4087","4103","4104","4105","4110","4114","4172","4174",
    "4176","4177","4194","4197","4225","4227","4237","4238","4273","4274",
    "4292","4294","4432","4433","4458","4459","4460","4472","4473","4476"
)
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $odiBowling.Cells.Item($i + 2, 2).Value = $bowlingCodes[$i]
}

# ---------------------------------------------------------------------------
# 4. "ODI Batting Extra" sheet - new, appended after "ODI Bowling".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$battingExtra.Cells.Item(1,2).Value = "BATTING_POSITION"
$battingExtra.Cells.Item(1,3).Value = "NUM_4"
$battingExtra.Cells.Item(1,4).Value = "NUM_6"
$battingExtra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4177","7","","","","NO"),
    @("4194","8","","","","NO"),
    @("4197","8","2","0","12.96%","NO"),
    @("4225","","","","","NO"),
    @("4227","","","","","NO"),
    @("4237","9","","","","NO"),
    @("4238","9","0","0","","NO"),
    @("4273","6","4","0","10.00%","NO"),
    @("4274","7","1","1","4.93%","NO"),
    @("4287","","","","","NO"),
    @("4292","8","0","0","0.83%","NO"),
    @("4294","","","","","NO"),
    @("4432","","","","","NO"),
    @("4433","","","","","NO"),
    @("4458","8","0","0","1.82%","NO"),
    @("4459","8","1","0","3.40%","NO"),
    @("4460","6","0","0","0.31%","NO"),
    @("4472","8","0","0","3.55%","NO"),
    @("4473","8","0","0","0.51%","NO"),
    @("4476","7","1","1","3.02%","NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $row = $extraRows[$i]
    $r = $i + 2
    $battingExtra.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $battingExtra.Cells.Item($r, 2).Value = [double]$row[1]
    }
    if ($row[2] -ne "") {
        $battingExtra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($row[3] -ne "") {
        $battingExtra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($row[4] -ne "") {
        $battingExtra.Cells.Item($r, 5).Value = $row[4]
    }
    $battingExtra.Cells.Item($r, 6).Value = $row[5]
}
